# Sheet Name, ListObject Name에 #으로 시작할 경우 제외
#
# The "Param_1 설명" helper column in the Character tables is itself named
# "#Param_1 설명" (already excluded by the generator because it starts with
# "#"). This change additionally marks the column's *type* row with a
# leading "#" (so "int32_cs" -> "#int32_cs") for the two in-sheet sample
# tables (Character_1 / Character_2), and fills in the previously-empty
# type/note cell for the second sample table's data row, so every row is
# consistently marked as excluded.

$wb = $excel.ActiveWorkbook
$sheetEnum = $wb.Worksheets.Item(1)
$sheetCharacter = $wb.Worksheets.Item(2)

# --- Data edits on the "Character" sheet -----------------------------------

# Table "Character_1" (rows 2-5): the Param_1 설명 column's type-row value
# becomes "#int32_cs" instead of "int32_cs".
$sheetCharacter.Range("J3").Value = "#int32_cs"

# Table "Character_2" (rows 9-11): same type-row change ...
$sheetCharacter.Range("J10").Value = "#int32_cs"

# ... and the previously-blank data cell for that column now carries the
# same "excluded" note the first table already had in row 4, with the
# leading apostrophe so Excel stores it as literal (quote-prefixed) text.
$sheetCharacter.Range("J11").Value = "'#으로 시작하는 컬럼은 제외"

# --- View-state bookkeeping --------------------------------------------------

# Move the UI selection on "Character" to J14, then hand the active tab
# back to "Enum" (tab 1) instead of leaving "#Character 설명" selected.
$sheetCharacter.Activate()
$sheetCharacter.Range("J14").Select()
$sheetEnum.Activate()
